$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray empty inline-string cells from row 18 (keep only the
# cells that actually carry content: A-D, Q, S).
$ws.Range("E18:P18").ClearContents()
$ws.Range("R18").ClearContents()
$ws.Range("T18:V18").ClearContents()

# Insert a new row above the current row 39 ("mental health intervention
# content") to hold a new ontology entry, shifting rows 39-66 down to 40-67.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the "measurement unit label" entry.
$ws.Range("A39").Value = "IAO:0000003"
$ws.Range("B39").Value = "measurement unit label"
$ws.Range("C39").Value = "A measurement unit label is as a label that is part of a scalar measurement datum and denotes a unit of measure."
$ws.Range("D39").Value = "datum label"
$ws.Range("P39").Value = "LSR1; LSR2; LSR3"
$ws.Range("Q39").Value = "Intervention content and delivery"
$ws.Range("S39").Value = "External"
